# update codes for gender, father's name
#
# - "Photo" column renamed to "Gender"
# - Row 2 "Name" changed from "Naw Nandar Oo" to "Nan Dae Thayi Sandar Aye"
# - The two photo-URL hyperlinks in column E are replaced with plain
#   "Male" / "Female" text (underlined, black) and their hyperlinks removed
# - Selection moved to F3

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Content updates -------------------------------------------------

# Row 2 "Name": Naw Nandar Oo -> Nan Dae Thayi Sandar Aye
$ws.Range("C2").Value = "Nan Dae Thayi Sandar Aye"

# Header: Photo -> Gender (also renames the Table1 column, since it is
# driven by the header cell text)
$ws.Range("E1").Value = "Gender"

# Replace the Photo-URL values with Male / Female, keeping the
# underlined look but switching off the "Hyperlink" blue color
# (automatic/black, theme color 1) to match a plain (non-link) style.
$ws.Range("E2").Value = "Male"
$ws.Range("E2").Font.ThemeColor = 1

$ws.Range("E3").Value = "Female"
$ws.Range("E3").Font.ThemeColor = 1

# --- Hyperlinks --------------------------------------------------------
# This engine's Hyperlinks.Delete() clears every hyperlink on the sheet
# (it is not scoped to the calling Range), so remove them all and
# re-create only the two mailto: links that must remain (B2 and B3).
$ws.Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("B2"), "mailto:sawmyintwin@gmail.com")
$ws.Hyperlinks.Add($ws.Range("B3"), "mailto:sawmyintwin@gmail.com")
# Re-assert the built-in Hyperlink cell style so B2/B3 keep their
# original look (re-adding the link can otherwise nudge the style).
$ws.Range("B2").Style = "Hyperlink"
$ws.Range("B3").Style = "Hyperlink"

# --- View ---------------------------------------------------------------
$ws.Range("F3").Select()
